$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain (non-ambiguous) text -- assign directly.
$ws.Range("D2").Value = "34.821.60"
$ws.Range("E2").Value = "  -2.45%  "
$ws.Range("D3").Value = "1.801.61"
$ws.Range("E3").Value = "  -3.18%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  -7.12%  "
$ws.Range("E9").Value = "  +3.11%  "
$ws.Range("E10").Value = "  -2.88%  "
$ws.Range("E11").Value = "  -1.50%  "
$ws.Range("D12").Value = "2.058.99"
$ws.Range("E12").Value = "  -3.48%  "
$ws.Range("D13").Value = "1.804.10"
$ws.Range("E13").Value = "  -3.52%  "
$ws.Range("E14").Value = "  -3.66%  "
$ws.Range("E15").Value = "  -6.41%  "
$ws.Range("E16").Value = "  -5.09%  "
$ws.Range("D17").Value = "34.707.82"
$ws.Range("E17").Value = "  -2.63%  "
$ws.Range("E18").Value = "  -2.19%  "
$ws.Range("D19").Value = "0.0₃0776"
$ws.Range("E19").Value = "  -3.22%  "
$ws.Range("E20").Value = "  -3.49%  "
$ws.Range("E21").Value = "  -4.55%  "
$ws.Range("E22").Value = "  -3.66%  "
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("E24").Value = "  -0.81%  "
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("E26").Value = "  -3.06%  "
$ws.Range("E27").Value = "  -3.96%  "
$ws.Range("E28").Value = "  -3.57%  "
$ws.Range("E29").Value = "  +8.24%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  +0.95%  "
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("E33").Value = "  -3.95%  "
$ws.Range("E34").Value = "  -8.12%  "
$ws.Range("E35").Value = "  +2.88%  "
$ws.Range("E36").Value = "  -2.81%  "
$ws.Range("E37").Value = "  +1.21%  "
$ws.Range("E38").Value = "  +0.52%  "
$ws.Range("D39").Value = "1.304.91"
$ws.Range("E39").Value = "  -2.93%  "
$ws.Range("E40").Value = "  -3.58%  "
$ws.Range("E41").Value = "  -0.34%  "
$ws.Range("E42").Value = "  -7.32%  "
$ws.Range("E43").Value = "  -5.38%  "
$ws.Range("E44").Value = "  -3.03%  "
$ws.Range("E45").Value = "  -13.32%  "
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("E47").Value = "  -2.88%  "
$ws.Range("D48").Value = "1.975.88"
$ws.Range("E49").Value = "  +0.19%  "
$ws.Range("E50").Value = "  +6.45%  "
$ws.Range("E51").Value = "  -6.03%  "

# Cells whose new values look like plain numbers (e.g. "39.00") would be
# auto-converted to numeric values by Excel. Force them to remain text by
# temporarily applying a text number format, then restore the default style
# so the resulting cell has no explicit style, matching the source data.
$numericLikeRefs = @("D5", "D6", "D8", "D9", "D10", "D14", "D15", "D16", "D20", "D21", "D22", "D24", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D36", "D37", "D38", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D50", "D51")
foreach ($ref in $numericLikeRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D5").Value = "231.45"
$ws.Range("D6").Value = "0.601"
$ws.Range("D8").Value = "39.00"
$ws.Range("D9").Value = "0.322"
$ws.Range("D10").Value = "0.0676"
$ws.Range("D14").Value = "0.654"
$ws.Range("D15").Value = "10.77"
$ws.Range("D16").Value = "4.53"
$ws.Range("D20").Value = "237.81"
$ws.Range("D21").Value = "11.63"
$ws.Range("D22").Value = "4.60"
$ws.Range("D24").Value = "2.20"
$ws.Range("D25").Value = "170.92"
$ws.Range("D26").Value = "7.68"
$ws.Range("D27").Value = "17.16"
$ws.Range("D29").Value = "1.55"
$ws.Range("D30").Value = "1.01"
$ws.Range("D31").Value = "3.97"
$ws.Range("D32").Value = "0.0546"
$ws.Range("D33").Value = "3.89"
$ws.Range("D34").Value = "1.75"
$ws.Range("D36").Value = "0.672"
$ws.Range("D37").Value = "89.79"
$ws.Range("D38").Value = "1.30"
$ws.Range("D40").Value = "0.0189"
$ws.Range("D41").Value = "2.46"
$ws.Range("D42").Value = "0.952"
$ws.Range("D43").Value = "14.31"
$ws.Range("D45").Value = "2.17"
$ws.Range("D46").Value = "6.11"
$ws.Range("D47").Value = "0.0507"
$ws.Range("D50").Value = "0.0663"
$ws.Range("D51").Value = "98.44"

foreach ($ref in $numericLikeRefs) {
    $ws.Range($ref).Style = "Normal"
}
